# Insert a new data row at row 232 (pushing existing rows 232-291 down to 233-292)
# and populate it with a new "Pomelo" price record for Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 232, shifting rows 232:291 down to 233:292
$ws.Rows.Item(232).Insert()

# Fill in the new row 232 with the new record's values
$ws.Cells.Item(232, 1).Value = 4
$ws.Cells.Item(232, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(232, 3).Value = "Los Lagos"
$ws.Cells.Item(232, 4).Value = 44711
$ws.Cells.Item(232, 5).Value = 10
$ws.Cells.Item(232, 6).Value = "Fruta"
$ws.Cells.Item(232, 7).Value = 100102
$ws.Cells.Item(232, 8).Value = "Cítricos"
$ws.Cells.Item(232, 9).Value = 100102006
$ws.Cells.Item(232, 10).Value = "Pomelo"
$ws.Cells.Item(232, 11).Value = "Start Ruby"
$ws.Cells.Item(232, 12).Value = "Primera"
$ws.Cells.Item(232, 13).Value = 80
$ws.Cells.Item(232, 14).Value = 14000
$ws.Cells.Item(232, 15).Value = 14000
$ws.Cells.Item(232, 16).Value = 14000
$ws.Cells.Item(232, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(232, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(232, 19).Value = 1000
$ws.Cells.Item(232, 20).Value = 14
